$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.186.92"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "2.525.74"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "323.94"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").Value = "109.12"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "0.529"
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.556"
$ws.Range("E9").Value = "  +3.86%  "
$ws.Range("D10").Value = "40.90"
$ws.Range("E10").Value = "  +5.26%  "
$ws.Range("D11").Value = "20.50"
$ws.Range("E11").Value = "  +12.01%  "
$ws.Range("D12").Value = "0.0826"
$ws.Range("E12").Value = "  +1.64%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").Value = "2.923.07"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "2.525.83"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("D18").Value = "48.024.48"
$ws.Range("E18").Value = "  +2.11%  "
$ws.Range("E19").Value = "  +4.26%  "
$ws.Range("D20").Value = "6.64"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "72.25"
$ws.Range("D24").Value = "269.56"
$ws.Range("E24").Value = "  +9.27%  "
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").Value = "26.22"
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").Value = "10.17"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("E29").Value = "  -3.53%  "
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").Value = "35.69"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("D32").Value = "49.60"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").Value = "19.95"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "0.0794"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("E37").Value = "  +1.46%  "
$ws.Range("E38").Value = "  +1.40%  "
$ws.Range("E39").Value = "  +1.71%  "
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("D41").Value = "22.20"
$ws.Range("E41").Value = "  +4.55%  "
$ws.Range("D42").Value = "119.99"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("E43").Value = "  -1.56%  "
$ws.Range("D44").Value = "0.0301"
$ws.Range("E44").Value = "  +1.79%  "
$ws.Range("D45").Value = "2.015.93"
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("D46").Value = "3.16"
$ws.Range("E46").Value = "  +4.20%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "1.90"
$ws.Range("E47").Value = "  +6.09%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "2.05"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("D50").Value = "5.27"
$ws.Range("E50").Value = "  +1.98%  "
$ws.Range("D51").Value = "79.48"
$ws.Range("E51").Value = "  +2.56%  "
